$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text (it holds numeric-looking
# strings like "132.60" that Excel would otherwise coerce into
# numbers and strip the trailing zero from).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '58.426.47'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.517.40'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '521.01'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = '132.60'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '0.558'
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").Value = '2.515.21'
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("D10").Value = '0.0975'
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = '5.16'
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '2.962.79'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").Value = '58.388.12'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '22.10'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '2.515.17'
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '321.10'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '4.15'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +8.59%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '64.78'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '0.405'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").Value = '7.37'
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("D29").Value = '0.0₃0754'
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  +2.19%  '
$ws.Range("D31").Value = '167.77'
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  +2.68%  '
$ws.Range("D33").Value = '6.25'
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").Value = '18.07'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  -6.78%  '
$ws.Range("D38").Value = '3.93'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("D40").Value = '36.25'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '0.770'
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("D42").Value = '277.93'
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = '3.48'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").Value = '4.98'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").Value = '127.63'
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("D47").Value = '0.0921'
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("D49").Value = '17.64'
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("D50").Value = '0.0213'
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").Value = '16.80'
$ws.Range("E51").Value = '  +0.58%  '

# Restore the default (unstyled) cell style on the Price column so
# the only change versus the original workbook is the cell content.
$ws.Range("D2:D51").Style = "Normal"
